# Refresh Price (column D) and Volume(1h) (column E) figures for the
# cryptos symbol list. Values are written with a leading apostrophe so
# Excel stores them as text (matching the sheet's existing inlineStr/text
# cells) instead of auto-converting to numbers/percentages; re-applying
# the "Normal" style afterwards clears the quote-prefix formatting that
# the apostrophe entry would otherwise leave behind.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'321.65"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "'-3.15%"
$ws.Range("E2").Style = "Normal"
$ws.Range("D3").Value = "'42.91"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "'-5.68%"
$ws.Range("E3").Style = "Normal"
$ws.Range("D4").Value = "'5.208"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "'-6.19%"
$ws.Range("E4").Style = "Normal"
$ws.Range("D5").Value = "'0.08177"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "'-2.30%"
$ws.Range("E5").Style = "Normal"
$ws.Range("E6").Value = "'-2.59%"
$ws.Range("E6").Style = "Normal"
$ws.Range("D7").Value = "'1.816"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "'-13.35%"
$ws.Range("E7").Style = "Normal"
$ws.Range("D8").Value = "'0.9352"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "'-5.55%"
$ws.Range("E8").Style = "Normal"
$ws.Range("D9").Value = "'0.1116"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "'-6.60%"
$ws.Range("E9").Style = "Normal"
$ws.Range("D10").Value = "'0.1858"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "'-3.79%"
$ws.Range("E10").Style = "Normal"
$ws.Range("D11").Value = "'0.09403"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "'-5.16%"
$ws.Range("E11").Style = "Normal"
$ws.Range("D12").Value = "'0.04624"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "'-1.08%"
$ws.Range("E12").Style = "Normal"
$ws.Range("D13").Value = "'7.449"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "'-22.09%"
$ws.Range("E13").Style = "Normal"
$ws.Range("D15").Value = "'0.001295"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "'0.42%"
$ws.Range("E15").Style = "Normal"
$ws.Range("D16").Value = "'0.005894"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "'-0.62%"
$ws.Range("E16").Style = "Normal"
$ws.Range("D17").Value = "'3.354"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "'-1.16%"
$ws.Range("E17").Style = "Normal"
$ws.Range("D18").Value = "'2.540"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "'-0.94%"
$ws.Range("E18").Style = "Normal"
$ws.Range("E19").Value = "'-0.65%"
$ws.Range("E19").Style = "Normal"
$ws.Range("D20").Value = "'0.1380"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "'1.93%"
$ws.Range("E20").Style = "Normal"
$ws.Range("E21").Value = "'-1.70%"
$ws.Range("E21").Style = "Normal"
$ws.Range("D22").Value = "'0.04139"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "'0.16%"
$ws.Range("E22").Style = "Normal"
$ws.Range("D23").Value = "'0.001244"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "'-3.79%"
$ws.Range("E23").Style = "Normal"
$ws.Range("D24").Value = "'0.004294"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "'-5.43%"
$ws.Range("E24").Style = "Normal"
$ws.Range("D25").Value = "'0.0001200"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "'-7.84%"
$ws.Range("E25").Style = "Normal"
$ws.Range("D26").Value = "'0.0002979"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "'-20.45%"
$ws.Range("E26").Style = "Normal"
$ws.Range("D38").Value = "'0.02727"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "'1.08%"
$ws.Range("E38").Style = "Normal"
$ws.Range("D39").Value = "'0.05561"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "'-3.54%"
$ws.Range("E39").Style = "Normal"
$ws.Range("D40").Value = "'0.008113"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "'2.77%"
$ws.Range("E40").Style = "Normal"
$ws.Range("E41").Value = "'-2.48%"
$ws.Range("E41").Style = "Normal"
$ws.Range("D42").Value = "'0.006541"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "'-17.68%"
$ws.Range("E42").Style = "Normal"
$ws.Range("D43").Value = "'0.002040"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "'0.82%"
$ws.Range("E43").Style = "Normal"
$ws.Range("D44").Value = "'0.008335"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "'-6.81%"
$ws.Range("E44").Style = "Normal"
$ws.Range("D45").Value = "'0.3506"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "'2.79%"
$ws.Range("E45").Style = "Normal"
$ws.Range("D46").Value = "'0.00006930"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "'-1.58%"
$ws.Range("E46").Style = "Normal"
$ws.Range("D47").Value = "'0.00000000750"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "'-0.14%"
$ws.Range("E47").Style = "Normal"
$ws.Range("D48").Value = "'0.003334"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "'-1.61%"
$ws.Range("E48").Style = "Normal"
$ws.Range("D49").Value = "'0.003531"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "'-0.16%"
$ws.Range("E49").Style = "Normal"
$ws.Range("D50").Value = "'0.00002100"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "'-0.14%"
$ws.Range("E50").Style = "Normal"
$ws.Range("D51").Value = "'0.0002000"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "'-0.14%"
$ws.Range("E51").Style = "Normal"
